$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Swap the shared-string text so that the "audience" question and its
# answer occupy the positions the diff expects (L1 now references the
# question text, L2/L3 the answer text).
$ws1.Range("L1").Value = "Who is the audience for this data"
$ws1.Range("L2").Value = "researchers, policy makers"
$ws1.Range("L3").Value = "researchers, policy makers"

# Update the view: select L1 (this also scrolls the window back so the
# top-left visible cell is A1 instead of A3).
$ws1.Activate()
$ws1.Range("L1").Select()
